# Insert a new price-report row for Alcachofa (Madrigal, Provincia del Elquí,
# 2022-08-04) above the existing row 26, pushing the subsequent rows (old
# 26:41) down to (27:42) and extending the sheet to A1:R42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 26; Excel shifts rows 26:41 down to 27:42 and grows
# the used range / dimension accordingly.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new record.
$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value = "Ñuble"
$ws.Cells.Item(26, 4).Value = 44777
$ws.Cells.Item(26, 5).Value = 16
$ws.Cells.Item(26, 6).Value = 100112013
$ws.Cells.Item(26, 7).Value = "Alcachofa"
$ws.Cells.Item(26, 8).Value = "Madrigal"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 60
$ws.Cells.Item(26, 11).Value = 14000
$ws.Cells.Item(26, 12).Value = 15000
$ws.Cells.Item(26, 13).Value = 14500
$ws.Cells.Item(26, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(26, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(26, 16).Value = 362
$ws.Cells.Item(26, 17).Value = 40
$ws.Cells.Item(26, 18).Value = "Hortaliza"
